$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update numeroDocumento (B2): 700100 -> 700101
$ws.Range("B2").Value = 700101

# Update usuario (D2): pruebasregistro48 -> pruebasregistro49
$ws.Range("D2").Value = "pruebasregistro49"

# Move the active cell / selection to D17 (cosmetic change in sheet view)
$ws.Range("D17").Select()
